$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows right after the header (before current row 2),
# shifting existing data rows down. This pushes the last 2 existing
# rows (21 and 22) out of the A1:C21 range, and row 22 is dropped
# entirely along with the overflow.
$ws.Rows.Item(2).Resize(2).Insert()

# Remove the now-extra trailing row so the used range stays at 22 rows
# worth of content shifted - the old row 22 (now row 24) and anything
# beyond the new dimension should be cleared.
$ws.Rows.Item(22).Resize(3).Delete()

# The insert copies the header row's bold/centered style onto the new
# rows; strip that back off so the new data rows stay unstyled like the
# rest of the data.
$ws.Rows.Item(2).Resize(2).ClearFormats()

# Populate the two newly inserted rows with the new accelerometer samples
$ws.Range("A2").Value = -0.03897037506103547
$ws.Range("B2").Value = -0.04966262578964246
$ws.Range("C2").Value = 0.6067050054669385

$ws.Range("A3").Value = -0.3774656057357791
$ws.Range("B3").Value = -0.08384630084037778
$ws.Range("C3").Value = 1.005563378334046
